# Insert a new weekly data row above row 143, shifting existing rows
# 143-212 down to 144-213 (dimension grows from A1:R212 to A1:R213).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 143; existing rows 143..212 shift to 144..213.
$ws.Rows.Item(143).Insert()

# Populate the newly inserted row 143 with the new weekly record.
$ws.Cells.Item(143, 1).Value = 4
$ws.Cells.Item(143, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(143, 3).Value = "Los Lagos"
$ws.Cells.Item(143, 4).Value = 44523
$ws.Cells.Item(143, 5).Value = 10
$ws.Cells.Item(143, 6).Value = 100114014
$ws.Cells.Item(143, 7).Value = "Betarraga"
$ws.Cells.Item(143, 8).Value = "Sin especificar"
$ws.Cells.Item(143, 9).Value = "Primera"
$ws.Cells.Item(143, 10).Value = 1200
$ws.Cells.Item(143, 11).Value = 900
$ws.Cells.Item(143, 12).Value = 1000
$ws.Cells.Item(143, 13).Value = 950
$ws.Cells.Item(143, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(143, 15).Value = "Región del Maule"
$ws.Cells.Item(143, 16).Value = 190
$ws.Cells.Item(143, 17).Value = 5
$ws.Cells.Item(143, 18).Value = "Hortaliza"
